$d = $word.ActiveDocument

# --- Paragraph 1: "On Pilgrimage - September 1955" (currently Heading1, wrapped by a bookmark) ---
$p1 = $d.Paragraphs.Item(1)

# Remove the bookmark that wraps this paragraph (on-pilgrimage---september-1955)
if ($d.Bookmarks.Exists("on-pilgrimage---september-1955")) {
    $d.Bookmarks.Item("on-pilgrimage---september-1955").Delete()
}

$r1 = $p1.Range
$r1.Text = "On Pilgrimage - September 1955"
$r1.Font.Bold = 0
$p1.Style = "Title"

# --- Paragraph 2: "By Dorothy Day" -> "Dorothy Day", style Authors, no bold ---
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.Text = "Dorothy Day"
$r2.Font.Bold = 0
$p2.Style = "Authors"
